$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Critical Minutes and Good Roaming Calculation (%)
$ws.Range("C3").Value = 2184
$ws.Range("D3").Value = 92.5

# Row 4 (Totals): Critical Minutes
$ws.Range("C4").Value = 2184

# Row 12: clear Driver Vintage date value
$ws.Range("E12").Value = ""

# Row 14: Total Samples
$ws.Range("B14").Value = 265400
